$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new trailing columns: ComponentAmount (F) and ComponentType (G),
# and append a new data row (row 29) for a bolt component on asset 2293.
$ws.Range("F1").Value = "ComponentAmount"

$ws.Range("A29").Value = 2293
$ws.Range("A29").NumberFormat = "0000"
$ws.Range("B29").Value = 10
$ws.Range("C29").Value = "Bolt"

$ws.Range("G1").Value = "ComponentType"
$ws.Range("G29").Value = "M24"

# Rename the "ComponentNumber" header in B1 to "ProductCode" (the column now
# tracks product codes instead of a component number, data unchanged).
$ws.Range("B1").Value = "ProductCode"

# Reflect the new filtering/overview behaviour: selecting B1 as the active
# cell (matches the clicked-header selection state captured on save).
$ws.Range("B1").Select()
